$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "human_health"

# Update row 2: PLA_impellar_virgin_(recipe_endpoint_h) -> PLA_virgin, value -> 0.409
$ws.Range("A2").Value = "PLA_virgin"
$ws.Range("B2").Value = 0.409

# Update row 3: PLA_impellar_recycled_(recipe_endpoint_h) -> PLA_recycled, value -> 0.164
$ws.Range("A3").Value = "PLA_recycled"
$ws.Range("B3").Value = 0.164

# Update row 4: PLA_impellar_recycled_industrial_(recipe_endpoint_h) -> PLA_recycled_industrial, value -> 0.32
$ws.Range("A4").Value = "PLA_recycled_industrial"
$ws.Range("B4").Value = 0.32
